$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 638.4
$ws.Range("I8").Value = 70
$ws.Range("J8").Value = 1491
$ws.Range("K8").Value = 210
$ws.Range("L8").Value = 4473
$ws.Range("M8").Value = -71
$ws.Range("N8").Value = -4751

$ws.Range("H52").Value = 3857.1428
$ws.Range("J52").Value = 3857.1428
$ws.Range("L52").Value = 11571.4284
$ws.Range("N52").Value = -11891.4284

$ws.Range("H88").Value = 3445.818
$ws.Range("I88").Value = 1300
$ws.Range("J88").Value = 3660.4
$ws.Range("K88").Value = 1300
$ws.Range("L88").Value = 3660.4
$ws.Range("M88").Value = -894
$ws.Range("N88").Value = -4472.4

$ws.Range("H91").Value = 3445.818
$ws.Range("I91").Value = 1300
$ws.Range("J91").Value = 3660.4
$ws.Range("K91").Value = 1300
$ws.Range("L91").Value = 3660.4
$ws.Range("M91").Value = 104
$ws.Range("N91").Value = -6468.4

$ws.Range("H101").Value = 619.4286
$ws.Range("I101").Value = 187.2
$ws.Range("J101").Value = 1700
$ws.Range("K101").Value = 561.5999999999999
$ws.Range("L101").Value = 5100
$ws.Range("M101").Value = 1060.4
$ws.Range("N101").Value = -8344

$ws.Range("H131").Value = 1533.9286
$ws.Range("I131").Value = 1622.9166
$ws.Range("K131").Value = 4868.7498
$ws.Range("M131").Value = 171.2502000000004

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H51").Value = 0
$ws.Range("J51").Value = 0
$ws.Range("L51").Value = 0
$ws.Range("N51").ClearContents()

$ws.Range("H74").Value = 4163.3145
$ws.Range("I74").Value = 671.5714
$ws.Range("K74").Value = 671.5714
$ws.Range("M74").Value = 202.4286

$ws.Range("H77").Value = 4163.3145
$ws.Range("I77").Value = 671.5714
$ws.Range("K77").Value = 3357.857
$ws.Range("M77").Value = 1010.143

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 476
$ws.Range("I8").Value = 476
$ws.Range("K8").Value = 476
$ws.Range("M8").Value = -336

$ws.Range("H64").Value = 1397.6786
$ws.Range("J64").Value = 746.0476
$ws.Range("L64").Value = 746.0476
$ws.Range("N64").Value = -1196.0476

$ws.Range("H67").Value = 1397.6786
$ws.Range("J67").Value = 746.0476
$ws.Range("L67").Value = 746.0476
$ws.Range("N67").Value = -2306.0476

$ws.Range("H86").Value = 2035.52
$ws.Range("I86").Value = 1804.1904
$ws.Range("J86").Value = 3250
$ws.Range("K86").Value = 1804.1904
$ws.Range("L86").Value = 3250
$ws.Range("M86").Value = -681.1904
$ws.Range("N86").Value = -5496

$ws.Range("H89").Value = 2035.52
$ws.Range("I89").Value = 1804.1904
$ws.Range("J89").Value = 3250
$ws.Range("K89").Value = 9020.951999999999
$ws.Range("L89").Value = 16250
$ws.Range("M89").Value = -3404.951999999999
$ws.Range("N89").Value = -27482

$ws.Range("H97").Value = 4426
$ws.Range("I97").Value = 4426
$ws.Range("K97").Value = 4426
$ws.Range("M97").Value = -3435

$ws.Range("H99").Value = 1468.5625
$ws.Range("I99").Value = 1445.2727
$ws.Range("J99").Value = 1519.8
$ws.Range("K99").Value = 1445.2727
$ws.Range("L99").Value = 1519.8
$ws.Range("M99").Value = 52.72730000000001
$ws.Range("N99").Value = -4515.8

$ws.Range("H111").Value = 35000
$ws.Range("J111").Value = 35000
$ws.Range("L111").Value = 35000
$ws.Range("N111").Value = -43180

$ws.Range("H134").Value = 30337.184
$ws.Range("I134").Value = 31076.027
$ws.Range("K134").Value = 93228.08099999999
$ws.Range("M134").Value = -90693.08099999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H25").Value = 1000
$ws.Range("I25").Value = 1000
$ws.Range("K25").Value = 1000
$ws.Range("M25").Value = -826

$ws.Range("H31").Value = 1218.5807
$ws.Range("I31").Value = 1034.2963
$ws.Range("J31").Value = 2462.5
$ws.Range("K31").Value = 1034.2963
$ws.Range("L31").Value = 2462.5
$ws.Range("M31").Value = -739.2963
$ws.Range("N31").Value = -3052.5

$ws.Range("H34").Value = 1218.5807
$ws.Range("I34").Value = 1034.2963
$ws.Range("J34").Value = 2462.5
$ws.Range("K34").Value = 1034.2963
$ws.Range("L34").Value = 2462.5
$ws.Range("M34").Value = -832.2963
$ws.Range("N34").Value = -2866.5

$ws.Range("H86").Value = 2180.3684
$ws.Range("I86").Value = 1685.9
$ws.Range("J86").Value = 2729.7778
$ws.Range("K86").Value = 1685.9
$ws.Range("L86").Value = 2729.7778
$ws.Range("M86").Value = -562.9000000000001
$ws.Range("N86").Value = -4975.7778

$ws.Range("H89").Value = 2180.3684
$ws.Range("I89").Value = 1685.9
$ws.Range("J89").Value = 2729.7778
$ws.Range("K89").Value = 8429.5
$ws.Range("L89").Value = 13648.889
$ws.Range("M89").Value = -2813.5
$ws.Range("N89").Value = -24880.889

$ws.Range("H134").Value = 5989.25
$ws.Range("I134").Value = 5989.25
$ws.Range("K134").Value = 17967.75
$ws.Range("M134").Value = -15432.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H126").Value = 5789.952
$ws.Range("I126").Value = 896.6667
$ws.Range("J126").Value = 6605.5
$ws.Range("K126").Value = 2690.0001
$ws.Range("L126").Value = 19816.5
$ws.Range("M126").Value = 2249.9999
$ws.Range("N126").Value = -29696.5

$ws.Range("H136").Value = 1765.7333
$ws.Range("I136").Value = 1457.1666
$ws.Range("K136").Value = 4371.4998
$ws.Range("M136").Value = 728.5002000000004

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 12000
$ws.Range("J58").Value = 12000
$ws.Range("L58").Value = 12000
$ws.Range("N58").Value = -12554

$ws.Range("H80").Value = 102507.586
$ws.Range("I80").Value = 3000
$ws.Range("J80").Value = 111553.73
$ws.Range("K80").Value = 3000
$ws.Range("L80").Value = 111553.73
$ws.Range("M80").Value = -2002
$ws.Range("N80").Value = -113549.73

$ws.Range("H83").Value = 102507.586
$ws.Range("I83").Value = 3000
$ws.Range("J83").Value = 111553.73
$ws.Range("K83").Value = 15000
$ws.Range("L83").Value = 557768.65
$ws.Range("M83").Value = -10008
$ws.Range("N83").Value = -567752.65

$ws.Range("H126").Value = 1903.9
$ws.Range("I126").Value = 1608
$ws.Range("J126").Value = 2199.8
$ws.Range("K126").Value = 4824
$ws.Range("L126").Value = 6599.400000000001
$ws.Range("M126").Value = -2354
$ws.Range("N126").Value = -11539.4

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 14998
$ws.Range("J39").Value = 14998
$ws.Range("L39").Value = 14998
$ws.Range("N39").Value = -15824

$ws.Range("H69").Value = 14333.333
$ws.Range("J69").Value = 14333.333
$ws.Range("L69").Value = 14333.333
$ws.Range("N69").Value = -15831.333

$ws.Range("H72").Value = 14333.333
$ws.Range("J72").Value = 14333.333
$ws.Range("L72").Value = 42999.999
$ws.Range("N72").Value = -50487.999
